$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 46 and Row 47 have had their match data swapped.
# Columns to swap: B, F, G, I, J, K, L, M, N, O, P, Q, U, V, W, X, AB
$cols = @("B","F","G","I","J","K","L","M","N","O","P","Q","U","V","W","X","AB")

foreach ($col in $cols) {
    $cell46 = $ws.Range($col + "46")
    $cell47 = $ws.Range($col + "47")
    $tmp = $cell46.Value2
    $cell46.Value2 = $cell47.Value2
    $cell47.Value2 = $tmp
}
